$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-21 14:46:52"
$zhcn.Range("H2").Value = "2016-03-21 14:47:14"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-21 14:46:56"
$dede.Range("H2").Value = "2016-03-21 14:47:20"
